# Applies the "blocks_pilot" data-table edits described in the commit:
#   fixed crash error, decreased ITI and increased trials during
#   driving+lexical, increased training trials of lexical only training

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- training_lexical row (row 2): increase training trials of lexical-only training ---
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = "no"

# --- lexical_wo_driving_roboto row (row 3) ---
$ws.Range("G3").Value = "no"

# --- lexical_wo_driving_neuefrutigerworld row (row 4) ---
$ws.Range("G4").Value = "no"

# --- training_driving_lexical row (row 5): decreased ITI, increased trials ---
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 8

# --- full_task_roboto row (row 6): decreased ITI, increased trials ---
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 12

# --- full_task_neuefrutigerworld row (row 7): decreased ITI, increased trials ---
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 15
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 12

# --- view update: move selection from F5 to G7 (also clears the old topLeftCell scroll) ---
$ws.Activate()
$ws.Range("G7").Select()
